$wb = $excel.ActiveWorkbook

# --- 1. EpgScreen: drop the topLeftCell scroll position (leave selection at M18) ---
$wsEpg = $wb.Worksheets.Item("EpgScreen")
$wsEpg.Activate()
$wsEpg.Range("M18").Select()

# --- 2. screenTitles: add a new ZapList / zaplijst row (row 14), matching row 13's style ---
$wsTitles = $wb.Worksheets.Item("screenTitles")
$wsTitles.Activate()
$wsTitles.Range("A14").Value = "ZapList"
$wsTitles.Range("B14").Value = "zaplijst"
$wsTitles.Range("A13:B13").Copy()
$wsTitles.Range("A14:B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsTitles.Range("B1").Select()

# --- 3. parameters: drop tabSelected, move the live selection ---
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Activate()
$wsParams.Range("B15").Select()

# --- 4. PIPScreen: resize column B, move the selection to A1:B1 ---
$wsPIP = $wb.Worksheets.Item("PIPScreen")
$wsPIP.Activate()
$wsPIP.Columns.Item(2).ColumnWidth = 8.8
$wsPIP.Range("A1:B1").Select()

# --- 5. Add the new DTVChannel sheet after PIPScreen ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "DTVChannel"

# Copy the PIPScreen header formatting (bold/centred/filled) onto row 1
$wsPIP.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill the cells in the same order the shared strings were first introduced
$newSheet.Range("A1").Value = "objectID"
$newSheet.Range("A2").Value = "UnAvailiabeChannelNumber"
$newSheet.Range("B2").Value = 321
$newSheet.Range("B1").Value = "ChannelNumber"

$newSheet.Range("A3").Value = "ExpectedFocousChannel"
$newSheet.Range("B3").Value = 323
$newSheet.Range("A4").Value = "FirstChannelNumber"
$newSheet.Range("B4").Value = 1
$newSheet.Range("A5").Value = "LastChannelNumber"
$newSheet.Range("B5").Value = 999
$newSheet.Range("A6").Value = "HDChannelNumber"
$newSheet.Range("B6").Value = 7
$newSheet.Range("A7").Value = "SDChannelNumber"
$newSheet.Range("B7").Value = 6

$newSheet.Columns.Item(1).ColumnWidth = 26.85546875
$newSheet.Columns.Item(2).ColumnWidth = 15.7109375

$newSheet.Range("B6").Select()
